$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.811.48"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.374.62"
$ws.Range("E3").Value = "  -5.50%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "474.40"
$ws.Range("E5").Value = "  -3.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.08"
$ws.Range("E6").Value = "  +1.46%  "

$ws.Range("E7").Value = "  +0.49%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.495"
$ws.Range("E8").Value = "  -3.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.376.29"
$ws.Range("E9").Value = "  -6.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0964"
$ws.Range("E10").Value = "  -2.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.49"
$ws.Range("E11").Value = "  -2.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.322"
$ws.Range("E12").Value = "  -3.43%  "

$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.802.24"
$ws.Range("E14").Value = "  -4.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.213.52"
$ws.Range("E15").Value = "  +0.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.20"
$ws.Range("E16").Value = "  -3.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  -3.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.390.41"
$ws.Range("E18").Value = "  -5.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.42"
$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "310.98"
$ws.Range("E20").Value = "  -3.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.68"
$ws.Range("E21").Value = "  -5.30%  "

$ws.Range("E22").Value = "  -0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.63"
$ws.Range("E23").Value = "  -3.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "56.70"
$ws.Range("E24").Value = "  -3.00%  "

$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.392"
$ws.Range("E26").Value = "  -4.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.157"
$ws.Range("E27").Value = "  -6.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.531.25"
$ws.Range("E28").Value = "  -3.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.23"
$ws.Range("E29").Value = "  -3.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0764"
$ws.Range("E30").Value = "  -3.29%  "

$ws.Range("E31").Value = "  +0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.90"
$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.80"
$ws.Range("E33").Value = "  -3.59%  "

$ws.Range("E34").Value = "  -2.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.97"
$ws.Range("E35").Value = "  -5.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.10"
$ws.Range("E36").Value = "  -4.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.842"
$ws.Range("E37").Value = "  -2.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.57"
$ws.Range("E38").Value = "  -3.59%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.35"
$ws.Range("E40").Value = "  -2.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.34"
$ws.Range("E41").Value = "  +1.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0541"
$ws.Range("E42").Value = "  -3.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.36"
$ws.Range("E43").Value = "  -4.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.582"
$ws.Range("E44").Value = "  -5.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0936"
$ws.Range("E45").Value = "  +3.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.21"
$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "253.66"
$ws.Range("E47").Value = "  -5.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.53"
$ws.Range("E48").Value = "  -5.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0220"
$ws.Range("E49").Value = "  -2.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.97"
$ws.Range("E50").Value = "  -5.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.789.01"
$ws.Range("E51").Value = "  -8.45%  "
